$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.908.05'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.38'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.524'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.63'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.262'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.871.84'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.653.56'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.575'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.15'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.899.85'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.01'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.59'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.92'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.37'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.15%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.72'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.415.92'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.37%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.892'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.555'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.917'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.12'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.20'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.41'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.97%  '
